# Regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# Update column G ("K") values for rows 2-38 on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 5
    3  = 7
    4  = 7
    5  = 4
    6  = 6
    7  = 8
    8  = 8
    9  = 7
    10 = 6
    11 = 9
    12 = 6
    13 = 6
    14 = 9
    15 = 5
    16 = 7
    17 = 11
    18 = 9
    19 = 9
    20 = 6
    21 = 10
    22 = 9
    23 = 5
    24 = 7
    25 = 5
    26 = 4
    27 = 7
    28 = 5
    29 = 4
    30 = 3
    31 = 5
    32 = 7
    33 = 8
    34 = 2
    35 = 3
    36 = 4
    37 = 3
    38 = 6
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
